$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tc_009")
$ws.Activate()

$ws.Range("A1").Value = "a2l_file_name/comment"
$ws.Range("I7").Select()
